$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Structural change: insert a new column before "desc" (old col B) so the
#    old "servizio" column (A) can be split into "id" (A) + "service" (B).
#    Doing the insert this way lets Excel naturally re-home the existing
#    column-width metadata (old col B width -> new col C) and shift the
#    existing per-column cell styles one column to the right, matching how
#    the real edit was produced.
# ---------------------------------------------------------------------------
$ws.Columns("B").Insert()

# ---------------------------------------------------------------------------
# 2. Capture "donor" cells for each distinct format we need, using cells
#    that already carry that exact formatting right after the insert above.
# ---------------------------------------------------------------------------
# style "bold + center"                (header, id/service)
$donorHeaderBold = $ws.Range("A1")
# style "bold + Arial/theme1 + center" (header, rest of columns)
$donorHeaderArial = $ws.Range("C1")
# style "plain + no alignment"         (id column, base font)
$donorPlain = $ws.Range("A2")
# style "plain + left align"           (not directly needed, but matches old desc col)
$donorPlainLeft = $ws.Range("C2")
# style "plain + center align"         (==> target style idx 7)
$donorPlainCenter = $ws.Range("D2")
# style "Arial/theme1 + no alignment"  (==> target style idx 4)
$donorArial = $ws.Range("I2")

# ---------------------------------------------------------------------------
# 3. Apply formats to the final target layout (A..I, rows 1..6) so every
#    cell ends up on the right style before we touch cell values.
# ---------------------------------------------------------------------------

# Header row
$donorHeaderBold.Copy()
$ws.Range("A1,B1").PasteSpecial(-4122)

$donorHeaderArial.Copy()
$ws.Range("C1,D1,E1,F1,G1,H1,I1").PasteSpecial(-4122)

# "id" column + the one leftover "Spotify Family" cell that keeps the plain font
$donorPlain.Copy()
$ws.Range("A2,A3,A4,A5,B5,A6").PasteSpecial(-4122)

# "service" column (minus B5) + price columns (G,H,I) -> Arial/theme1, no alignment
$donorArial.Copy()
$ws.Range("B2,B3,B4,B6,G2,H2,I2,G3,H3,I3,G4,H4,I4,G5,H5,I5,G6,H6,I6").PasteSpecial(-4122)

# "desc" column -> Arial/theme1 font, left aligned (derive: start from Arial/no-align, then set alignment)
$donorArial.Copy()
$ws.Range("C2,C3,C4,C5,C6").PasteSpecial(-4122)
$ws.Range("C2,C3,C4,C5,C6").HorizontalAlignment = -4131

# "multiaccount" (row 2 only), "account", "custom" -> Arial/theme1 font, centered
$donorArial.Copy()
$ws.Range("D2,E2,F2,E3,F3,E4,F4,E5,F5,E6,F6").PasteSpecial(-4122)
$ws.Range("D2,E2,F2,E3,F3,E4,F4,E5,F5,E6,F6").HorizontalAlignment = -4108

# "multiaccount" rows 3-6 -> plain font, centered (matches target exactly)
$donorPlainCenter.Copy()
$ws.Range("D3,D4,D5,D6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Column width: "desc" (now column C) keeps the 29-wide custom width.
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 28.166667

# ---------------------------------------------------------------------------
# 5. Cell values (header row)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "service"
$ws.Range("C1").Value = "desc"
$ws.Range("D1").Value = "multiaccount"
$ws.Range("E1").Value = "account"
$ws.Range("F1").Value = "custom"
$ws.Range("G1").Value = "it"
$ws.Range("H1").Value = "es"
$ws.Range("I1").Value = "gb"

# ---------------------------------------------------------------------------
# 6. Cell values (data rows). "A" becomes a machine id/slug, "B" keeps the
#    human readable service name that used to live in the old "A" column.
# ---------------------------------------------------------------------------
# Row 2 - Netflix
$ws.Range("B2").Value = "Netflix"
$ws.Range("A2").Value = "netflix"
$ws.Range("C2").Value = "Servizio streaming video Ultra HD"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 17.99
$ws.Range("H2").Value = 20.99
$ws.Range("I2").Value = 14.99

# Row 3 - Amazon Music
$ws.Range("B3").Value = "Amazon Music"
$ws.Range("A3").Value = "amazonmusic"
$ws.Range("C3").Value = "Servizio musica streaming "
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 13.99
$ws.Range("H3").Value = 17.99
$ws.Range("I3").Value = 16.99

# Row 4 - Apple Music
$ws.Range("B4").Value = "Apple Music"
$ws.Range("A4").Value = "applemusic"
$ws.Range("C4").Value = "Servizio musica streaming"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 15.99
$ws.Range("H4").Value = 18.99
$ws.Range("I4").Value = 15.99

# Row 5 - Spotify Family
$ws.Range("B5").Value = "Spotify Family"
$ws.Range("A5").Value = "spotifyfamily"
$ws.Range("C5").Value = "Servizio musica streaming"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 14.99
$ws.Range("H5").Value = 16.99
$ws.Range("I5").Value = 12.99

# Row 6 - Custom
$ws.Range("B6").Value = "Custom"
$ws.Range("A6").Value = "custom"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
